# Updated daily log and folder structure v3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New log entries (row 7 and row 8): dates, activities, shared strings
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 43922
$ws.Range("B7").Value = "Getting started with UPPMAX, logging in and creating soft linkes. Learning how to find things in the directories."

$ws.Range("A8").Value = 43924
$ws.Range("B8").Value = "Running my first analysis of the data using FastQC. Learning how to connect to a working node and running in interactive mode. Synchronizing everything with my git-repository. Creating gitignore-file to prevent adding to large files"
$ws.Range("C8").Value = "Analyze result from fastQC and then do the DNA assembly"

# ---------------------------------------------------------------------------
# 2) Row heights for the data rows
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(3).RowHeight = 52
$ws.Rows.Item(4).RowHeight = 64
$ws.Rows.Item(5).RowHeight = 35
$ws.Rows.Item(6).RowHeight = 29
$ws.Rows.Item(7).RowHeight = 55
$ws.Rows.Item(8).RowHeight = 114

# ---------------------------------------------------------------------------
# 3) Column widths for B (What did I do) and C (Future plans)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 39
$ws.Columns.Item(3).ColumnWidth = 32.666666666666664

# ---------------------------------------------------------------------------
# 4) Fonts / colors: move header + body text from theme colors to explicit
#    black, and the header fill from a theme tint to an explicit light blue.
# ---------------------------------------------------------------------------
# Header row
$ws.Range("A1:D1").Font.Color = 0
$ws.Range("A1:D1").Interior.Color = 16247773
$ws.Range("A1:D1").Interior.PatternColor = 0

# Date column
$ws.Range("A2:A8").Font.Color = 0

# Remaining body / blank cells
$ws.Range("B2:D28").Font.Color = 0
$ws.Range("A9:A28").Font.Color = 0

# ---------------------------------------------------------------------------
# 5) Selection
# ---------------------------------------------------------------------------
$sel = $ws.Range("B4").Select()
